$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): Right count 5 -> 4, Wrong count -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 (Total): Right total 55 -> 44, Wrong total -5 -> -10
$ws.Range("B12").Value = 44
$ws.Range("C12").Value = -10

# E12 score summary text updated to match new totals
$ws.Range("E12").Value = "34 / 112"
